$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 24220.857
$ws.Range("J87").Value = 24220.857
$ws.Range("L87").Value = 24220.857
$ws.Range("N87").Value = -26716.857
$ws.Range("H90").Value = 24220.857
$ws.Range("J90").Value = 24220.857
$ws.Range("L90").Value = 72662.571
$ws.Range("N90").Value = -85142.571
$ws.Range("H117").Value = 48667.832
$ws.Range("J117").Value = 48667.832
$ws.Range("L117").Value = 48667.832
$ws.Range("N117").Value = -57845.832
$ws.Range("H133").Value = 46177.668
$ws.Range("J133").Value = 46177.668
$ws.Range("L133").Value = 46177.668
$ws.Range("N133").Value = -56297.668

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 38315
$ws.Range("J80").Value = 38315
$ws.Range("L80").Value = 38315
$ws.Range("N80").Value = -40311
$ws.Range("H83").Value = 38315
$ws.Range("J83").Value = 38315
$ws.Range("L83").Value = 114945
$ws.Range("N83").Value = -124929
$ws.Range("H114").Value = 45948
$ws.Range("J114").Value = 45948
$ws.Range("L114").Value = 45948
$ws.Range("N114").Value = -54626

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 55219.5
$ws.Range("J57").Value = 55219.5
$ws.Range("L57").Value = 55219.5
$ws.Range("N57").Value = -56659.5
$ws.Range("H100").Value = 42656
$ws.Range("J100").Value = 42656
$ws.Range("L100").Value = 42656
$ws.Range("N100").Value = -44820
$ws.Range("H102").Value = 18301.2
$ws.Range("I102").Value = 11506
$ws.Range("K102").Value = 11506
$ws.Range("M102").Value = -8261
$ws.Range("H103").Value = 42308.2
$ws.Range("J103").Value = 42308.2
$ws.Range("L103").Value = 42308.2
$ws.Range("N103").Value = -44652.2
$ws.Range("H122").Value = 40132.668
$ws.Range("J122").Value = 40132.668
$ws.Range("L122").Value = 40132.668
$ws.Range("N122").Value = -49932.668
$ws.Range("H133").Value = 48950
$ws.Range("J133").Value = 48950
$ws.Range("L133").Value = 48950
$ws.Range("N133").Value = -59070
$ws.Range("H136").Value = 55219.5
$ws.Range("J136").Value = 55219.5
$ws.Range("L136").Value = 55219.5
$ws.Range("N136").Value = -65419.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H52").Value = 66533
$ws.Range("J52").Value = 66533
$ws.Range("L52").Value = 66533
$ws.Range("N52").Value = -67121
$ws.Range("H82").Value = 44181
$ws.Range("J82").Value = 44181
$ws.Range("L82").Value = 44181
$ws.Range("N82").Value = -44903
$ws.Range("H85").Value = 44181
$ws.Range("J85").Value = 44181
$ws.Range("L85").Value = 44181
$ws.Range("N85").Value = -46677
$ws.Range("H88").Value = 45339
$ws.Range("J88").Value = 45339
$ws.Range("L88").Value = 45339
$ws.Range("N88").Value = -46151
$ws.Range("H91").Value = 45339
$ws.Range("J91").Value = 45339
$ws.Range("L91").Value = 45339
$ws.Range("N91").Value = -48147
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H124").Value = 37248
$ws.Range("J124").Value = 37248
$ws.Range("L124").Value = 37248
$ws.Range("N124").Value = -42158
$ws.Range("H125").Value = 29658.5
$ws.Range("J125").Value = 29658.5
$ws.Range("L125").Value = 29658.5
$ws.Range("N125").Value = -34578.5
$ws.Range("H131").Value = 34196
$ws.Range("J131").Value = 34196
$ws.Range("L131").Value = 34196
$ws.Range("N131").Value = -44276
$ws.Range("H137").Value = 35384.617
$ws.Range("J137").Value = 35384.617
$ws.Range("L137").Value = 35384.617
$ws.Range("N137").Value = -45584.617
$ws.Range("H139").Value = 59679.8
$ws.Range("J139").Value = 64099.75
$ws.Range("L139").Value = 64099.75
$ws.Range("N139").Value = -74379.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 9758.666999999999
$ws.Range("I122").Value = 566.1111
$ws.Range("J122").Value = 37336.332
$ws.Range("K122").Value = 5094.9999
$ws.Range("L122").Value = 336026.988
$ws.Range("M122").Value = -2644.9999
$ws.Range("N122").Value = -340926.988

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 30067.5
$ws.Range("J86").Value = 30067.5
$ws.Range("L86").Value = 30067.5
$ws.Range("N86").Value = -32439.5
$ws.Range("H89").Value = 30067.5
$ws.Range("J89").Value = 30067.5
$ws.Range("L89").Value = 90202.5
$ws.Range("N89").Value = -102058.5
$ws.Range("H110").Value = 47683.668
$ws.Range("J110").Value = 47683.668
$ws.Range("L110").Value = 47683.668
$ws.Range("N110").Value = -55863.668
$ws.Range("H119").Value = 48566
$ws.Range("J119").Value = 48566
$ws.Range("L119").Value = 48566
$ws.Range("N119").Value = -58242
$ws.Range("H127").Value = 36316.2
$ws.Range("J127").Value = 36316.2
$ws.Range("L127").Value = 36316.2
$ws.Range("N127").Value = -46236.2
$ws.Range("H131").Value = 29330.666
$ws.Range("J131").Value = 29330.666
$ws.Range("L131").Value = 29330.666
$ws.Range("N131").Value = -39410.666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1465.6842
$ws.Range("I16").Value = 1371.0938
$ws.Range("J16").Value = 1970.1666
$ws.Range("K16").Value = 1371.0938
$ws.Range("L16").Value = 1970.1666
$ws.Range("M16").Value = -1201.0938
$ws.Range("N16").Value = -2310.1666
$ws.Range("H36").Value = 48640
$ws.Range("J36").Value = 48640
$ws.Range("L36").Value = 48640
$ws.Range("N36").Value = -49764
$ws.Range("H88").Value = 43077
$ws.Range("J88").Value = 43077
$ws.Range("L88").Value = 43077
$ws.Range("N88").Value = -43933
$ws.Range("H91").Value = 43077
$ws.Range("J91").Value = 43077
$ws.Range("L91").Value = 43077
$ws.Range("N91").Value = -46041
$ws.Range("H96").Value = 37594.5
$ws.Range("J96").Value = 37594.5
$ws.Range("L96").Value = 37594.5
$ws.Range("N96").Value = -43086.5
$ws.Range("H99").Value = 25411.428
$ws.Range("I99").Value = 12626.667
$ws.Range("K99").Value = 12626.667
$ws.Range("M99").Value = -9631.666999999999
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H123").Value = 38134.332
$ws.Range("J123").Value = 38134.332
$ws.Range("L123").Value = 38134.332
$ws.Range("N123").Value = -47934.332
$ws.Range("H124").Value = 48214
$ws.Range("J124").Value = 48214
$ws.Range("L124").Value = 48214
$ws.Range("N124").Value = -58034
$ws.Range("H129").Value = 45429
$ws.Range("J129").Value = 45429
$ws.Range("L129").Value = 45429
$ws.Range("N129").Value = -55429
$ws.Range("H131").Value = 44251
$ws.Range("J131").Value = 44251
$ws.Range("L131").Value = 44251
$ws.Range("N131").Value = -54331
$ws.Range("H133").Value = 39700
$ws.Range("J133").Value = 39700
$ws.Range("L133").Value = 39700
$ws.Range("N133").Value = -44760
$ws.Range("H137").Value = 40883.332
$ws.Range("J137").Value = 40883.332
$ws.Range("L137").Value = 40883.332
$ws.Range("N137").Value = -51083.332
$ws.Range("H139").Value = 50319.6
$ws.Range("J139").Value = 50319.6
$ws.Range("L139").Value = 50319.6
$ws.Range("N139").Value = -60599.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H103").Value = 39182
$ws.Range("J103").Value = 39182
$ws.Range("L103").Value = 39182
$ws.Range("N103").Value = -41526
$ws.Range("H106").Value = 29492.285
$ws.Range("J106").Value = 29492.285
$ws.Range("L106").Value = 29492.285
$ws.Range("N106").Value = -32016.285
$ws.Range("H139").Value = 51219.8
$ws.Range("J139").Value = 51219.8
$ws.Range("L139").Value = 51219.8
$ws.Range("N139").Value = -61499.8
